$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Delete()
[void]$ws.Rows.Item(1).Select()
